$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Add a date value in C1 (serial 44307 = 2021-04-21), formatted with the
# built-in short-date number format (numFmtId 14).
$ws.Range("C1").Value = 44307
$ws.Range("C1").NumberFormat = "mm-dd-yy"
